$wb = $excel.ActiveWorkbook

# The "trial_investigators" sheet had a bad/duplicate test row for
# "Mark Andrew Dickson" (row 3) - remove it. Excel drops the now-unused
# "Mark Andrew Dickson" shared string automatically and the rows below
# shift up.
$ws1 = $wb.Worksheets.Item("trial_investigators")
$ws1.Rows(3).Delete()

# Re-point the active sheet/selection at the corrected sheet, matching
# the saved view state after the edit.
$ws1.Activate()
$ws1.Range("A3").Select()
